$wb = $excel.ActiveWorkbook

# --- Add the new "top20" worksheet as the last sheet in the workbook ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$ws.Name = "top20"

# --- Populate the "Top 20" nodes / edges tables ---
$ws.Range("A1").Value2 = "Opencare"
$ws.Range("I1").Value2 = "AllObese"
$ws.Range("B2").Value2 = "degree"
$ws.Range("E2").Value2 = "betweenness"
$ws.Range("J2").Value2 = "degree"
$ws.Range("M2").Value2 = "betweenness"
$ws.Range("A3").Value2 = "research question"
$ws.Range("B3").Value2 = 353
$ws.Range("D3").Value2 = "design intervention"
$ws.Range("E3").Value2 = 47230.1769163
$ws.Range("I3").Value2 = "lakedistrict"
$ws.Range("J3").Value2 = 307
$ws.Range("L3").Value2 = "love"
$ws.Range("M3").Value2 = 194109.218623
$ws.Range("A4").Value2 = "community-based care"
$ws.Range("B4").Value2 = 294
$ws.Range("D4").Value2 = "research question"
$ws.Range("E4").Value2 = 45931.2935821
$ws.Range("I4").Value2 = "cumbria"
$ws.Range("J4").Value2 = 302
$ws.Range("L4").Value2 = "doncaster"
$ws.Range("M4").Value2 = 169907.925459
$ws.Range("A5").Value2 = "migration"
$ws.Range("B5").Value2 = 252
$ws.Range("D5").Value2 = "community-based care"
$ws.Range("E5").Value2 = 33484.0121866
$ws.Range("I5").Value2 = "nature"
$ws.Range("J5").Value2 = 296
$ws.Range("L5").Value2 = "yorkshire"
$ws.Range("M5").Value2 = 169234.078934
$ws.Range("A6").Value2 = "design intervention"
$ws.Range("B6").Value2 = 242
$ws.Range("D6").Value2 = "migration"
$ws.Range("E6").Value2 = 28559.9297199
$ws.Range("I6").Value2 = "photooftheday"
$ws.Range("J6").Value2 = 276
$ws.Range("L6").Value2 = "fitness"
$ws.Range("M6").Value2 = 167930.819144
$ws.Range("A7").Value2 = "legality"
$ws.Range("B7").Value2 = 207
$ws.Range("D7").Value2 = "participatory design/collaboration"
$ws.Range("E7").Value2 = 20259.3309153
$ws.Range("I7").Value2 = "england"
$ws.Range("J7").Value2 = 271
$ws.Range("L7").Value2 = "photooftheday"
$ws.Range("M7").Value2 = 152372.254934
$ws.Range("A8").Value2 = "existing system failure"
$ws.Range("B8").Value2 = 205
$ws.Range("D8").Value2 = "mental health"
$ws.Range("E8").Value2 = 18396.7449126
$ws.Range("I8").Value2 = "love"
$ws.Range("J8").Value2 = 270
$ws.Range("L8").Value2 = "england"
$ws.Range("M8").Value2 = 116517.73459
$ws.Range("A9").Value2 = "resource strain"
$ws.Range("B9").Value2 = 197
$ws.Range("D9").Value2 = "legality"
$ws.Range("E9").Value2 = 18201.6317347
$ws.Range("I9").Value2 = "doncaster"
$ws.Range("J9").Value2 = 248
$ws.Range("L9").Value2 = "nature"
$ws.Range("M9").Value2 = 116455.52327
$ws.Range("A10").Value2 = "autonomy"
$ws.Range("B10").Value2 = 194
$ws.Range("D10").Value2 = "sustainability"
$ws.Range("E10").Value2 = 17970.4238949
$ws.Range("I10").Value2 = "yorkshire"
$ws.Range("J10").Value2 = 233
$ws.Range("L10").Value2 = "beer"
$ws.Range("M10").Value2 = 115651.48626
$ws.Range("A11").Value2 = "story sharing"
$ws.Range("B11").Value2 = 192
$ws.Range("D11").Value2 = "OpenSource"
$ws.Range("E11").Value2 = 17185.3229219
$ws.Range("I11").Value2 = "fitness"
$ws.Range("J11").Value2 = 219
$ws.Range("L11").Value2 = "lakedistrict"
$ws.Range("M11").Value2 = 111755.090618
$ws.Range("A12").Value2 = "self-care"
$ws.Range("B12").Value2 = 187
$ws.Range("D12").Value2 = "resource strain"
$ws.Range("E12").Value2 = 16045.663783
$ws.Range("I12").Value2 = "instagood"
$ws.Range("J12").Value2 = 190
$ws.Range("L12").Value2 = "london"
$ws.Range("M12").Value2 = 106343.338891
$ws.Range("A13").Value2 = "mental health"
$ws.Range("B13").Value2 = 186
$ws.Range("D13").Value2 = "existing system failure"
$ws.Range("E13").Value2 = 15514.6669225
$ws.Range("I13").Value2 = "autumn"
$ws.Range("J13").Value2 = 185
$ws.Range("L13").Value2 = "cumbria"
$ws.Range("M13").Value2 = 101469.14643
$ws.Range("A14").Value2 = "participatory design/collaboration"
$ws.Range("B14").Value2 = 174
$ws.Range("D14").Value2 = "story sharing"
$ws.Range("E14").Value2 = 15269.9348831
$ws.Range("I14").Value2 = "somerset"
$ws.Range("J14").Value2 = 181
$ws.Range("L14").Value2 = "photography"
$ws.Range("M14").Value2 = 90807.2312481
$ws.Range("A15").Value2 = "OpenSource"
$ws.Range("B15").Value2 = 168
$ws.Range("D15").Value2 = "autonomy"
$ws.Range("E15").Value2 = 13828.0824899
$ws.Range("I15").Value2 = "uk"
$ws.Range("J15").Value2 = 176
$ws.Range("L15").Value2 = "summer"
$ws.Range("M15").Value2 = 81590.2690378
$ws.Range("A16").Value2 = "sustainability"
$ws.Range("B16").Value2 = 167
$ws.Range("D16").Value2 = "outside existing systems"
$ws.Range("E16").Value2 = 12485.1297653
$ws.Range("I16").Value2 = "summer"
$ws.Range("J16").Value2 = 175
$ws.Range("L16").Value2 = "somerset"
$ws.Range("M16").Value2 = 73216.1552298
$ws.Range("A17").Value2 = "outside existing systems"
$ws.Range("B17").Value2 = 153
$ws.Range("D17").Value2 = "youth"
$ws.Range("E17").Value2 = 9116.73355706
$ws.Range("I17").Value2 = "beautiful"
$ws.Range("J17").Value2 = 168
$ws.Range("L17").Value2 = "uk"
$ws.Range("M17").Value2 = 72370.9795224
$ws.Range("A18").Value2 = "skill sharing"
$ws.Range("B18").Value2 = 147
$ws.Range("D18").Value2 = "self-care"
$ws.Range("E18").Value2 = 9102.77860101
$ws.Range("I18").Value2 = "vscocam"
$ws.Range("J18").Value2 = 168
$ws.Range("L18").Value2 = "picoftheday"
$ws.Range("M18").Value2 = 69732.6299631
$ws.Range("A19").Value2 = "youth"
$ws.Range("B19").Value2 = 146
$ws.Range("D19").Value2 = "education"
$ws.Range("E19").Value2 = 8877.96368386
$ws.Range("I19").Value2 = "photography"
$ws.Range("J19").Value2 = 164
$ws.Range("L19").Value2 = "autumn"
$ws.Range("M19").Value2 = 65086.7064403
$ws.Range("A20").Value2 = "care networks"
$ws.Range("B20").Value2 = 144
$ws.Range("D20").Value2 = "care networks"
$ws.Range("E20").Value2 = 8467.05745516
$ws.Range("I20").Value2 = "bawtry"
$ws.Range("J20").Value2 = 149
$ws.Range("L20").Value2 = "instagood"
$ws.Range("M20").Value2 = 63939.8904461
$ws.Range("A21").Value2 = "meaningful life"
$ws.Range("B21").Value2 = 136
$ws.Range("D21").Value2 = "place-based"
$ws.Range("E21").Value2 = 8098.45290202
$ws.Range("I21").Value2 = "picoftheday"
$ws.Range("J21").Value2 = 145
$ws.Range("L21").Value2 = "glitter"
$ws.Range("M21").Value2 = 60498.0279691
$ws.Range("A22").Value2 = "intergenerational"
$ws.Range("B22").Value2 = 133
$ws.Range("D22").Value2 = "housing"
$ws.Range("E22").Value2 = 7688.34130805
$ws.Range("I22").Value2 = "landscape"
$ws.Range("J22").Value2 = 144
$ws.Range("L22").Value2 = "bawtry"
$ws.Range("M22").Value2 = 60462.5935553
$ws.Range("B25").Value2 = "edge_force"
$ws.Range("E25").Value2 = "betweenness"
$ws.Range("J25").Value2 = "weigth"
$ws.Range("M25").Value2 = "betweenness"
$ws.Range("A26").Value2 = "research question - community-based care"
$ws.Range("B26").Value2 = 12
$ws.Range("D26").Value2 = "conflict - resource strain"
$ws.Range("E26").Value2 = 2044.99936283
$ws.Range("I26").Value2 = "cumbria - lakedistrict"
$ws.Range("J26").Value2 = 277
$ws.Range("L26").Value2 = "yorkshire - beer"
$ws.Range("M26").Value2 = 43411.9606449
$ws.Range("A27").Value2 = "migration - research question"
$ws.Range("B27").Value2 = 11
$ws.Range("D27").Value2 = "grassroots - making rules for spaces"
$ws.Range("E27").Value2 = 1976
$ws.Range("I27").Value2 = "thegentlemansretreat - bawtry"
$ws.Range("J27").Value2 = 188
$ws.Range("L27").Value2 = "blue - glitter"
$ws.Range("M27").Value2 = 38398.1598062
$ws.Range("A28").Value2 = "migration - building relationships"
$ws.Range("B28").Value2 = 11
$ws.Range("D28").Value2 = "medical research - nutrition"
$ws.Range("E28").Value2 = 1976
$ws.Range("I28").Value2 = "bawtry - barbershop"
$ws.Range("J28").Value2 = 183
$ws.Range("L28").Value2 = "london - bluray"
$ws.Range("M28").Value2 = 24442.0897601
$ws.Range("A29").Value2 = "community-based care - legality"
$ws.Range("B29").Value2 = 10
$ws.Range("D29").Value2 = "legality - informal discussion"
$ws.Range("E29").Value2 = 1976
$ws.Range("I29").Value2 = "thegentlemansretreat - barbershop"
$ws.Range("J29").Value2 = 183
$ws.Range("L29").Value2 = "louth - nature"
$ws.Range("M29").Value2 = 23930.2236418
$ws.Range("A30").Value2 = "migration - resource strain"
$ws.Range("B30").Value2 = 10
$ws.Range("D30").Value2 = "design intervention - falling"
$ws.Range("E30").Value2 = 1377.65351504
$ws.Range("I30").Value2 = "bawtry - themanclub"
$ws.Range("J30").Value2 = 181
$ws.Range("L30").Value2 = "beer - view"
$ws.Range("M30").Value2 = 22979.44499
$ws.Range("A31").Value2 = "resource strain - research question"
$ws.Range("B31").Value2 = 10
$ws.Range("D31").Value2 = "wearable technology - design intervention"
$ws.Range("E31").Value2 = 1129.40809968
$ws.Range("I31").Value2 = "thegentlemansretreat - themanclub"
$ws.Range("J31").Value2 = 181
$ws.Range("L31").Value2 = "yorkshire - fitness"
$ws.Range("M31").Value2 = 22811.1542731
$ws.Range("A32").Value2 = "resource strain - community-based care"
$ws.Range("B32").Value2 = 9
$ws.Range("D32").Value2 = "migration - politics of healthcare"
$ws.Range("E32").Value2 = 1119.5141744
$ws.Range("I32").Value2 = "themanclub - barbershop"
$ws.Range("J32").Value2 = 177
$ws.Range("L32").Value2 = "beer - northyorkshire"
$ws.Range("M32").Value2 = 18415.3547227
$ws.Range("A33").Value2 = "legality - migration"
$ws.Range("B33").Value2 = 9
$ws.Range("D33").Value2 = "migration - design intervention"
$ws.Range("E33").Value2 = 1074.50991458
$ws.Range("I33").Value2 = "bawtry - apothecary87"
$ws.Range("J33").Value2 = 172
$ws.Range("L33").Value2 = "trainerlife - picoftheday"
$ws.Range("M33").Value2 = 18272.8558043
$ws.Range("A34").Value2 = "existing system failure - legality"
$ws.Range("B34").Value2 = 9
$ws.Range("D34").Value2 = "governance - legality"
$ws.Range("E34").Value2 = 1051.18912577
$ws.Range("I34").Value2 = "thegentlemansretreat - apothecary87"
$ws.Range("J34").Value2 = 172
$ws.Range("L34").Value2 = "bawtry - beer"
$ws.Range("M34").Value2 = 17954.6204214
$ws.Range("A35").Value2 = "migration - story sharing"
$ws.Range("B35").Value2 = 9
$ws.Range("D35").Value2 = "design intervention - research question"
$ws.Range("E35").Value2 = 1009.49255368
$ws.Range("I35").Value2 = "themanclub - apothecary87"
$ws.Range("J35").Value2 = 171
$ws.Range("L35").Value2 = "fitness - lakedistrict"
$ws.Range("M35").Value2 = 16643.1789459
$ws.Range("A36").Value2 = "legality - research question"
$ws.Range("B36").Value2 = 9
$ws.Range("D36").Value2 = "holistic healthcare - conceptual framework"
$ws.Range("E36").Value2 = 989
$ws.Range("I36").Value2 = "thegentlemansretreat - barberlife"
$ws.Range("J36").Value2 = 168
$ws.Range("L36").Value2 = "love - saturday"
$ws.Range("M36").Value2 = 16413.4765776
$ws.Range("A37").Value2 = "safety - regulation"
$ws.Range("B37").Value2 = 9
$ws.Range("D37").Value2 = "instagram - medical professionals"
$ws.Range("E37").Value2 = 989
$ws.Range("I37").Value2 = "barbershop - apothecary87"
$ws.Range("J37").Value2 = 168
$ws.Range("L37").Value2 = "yorkshirevapers - northyorkshire"
$ws.Range("M37").Value2 = 16296
$ws.Range("A38").Value2 = "mental health - creativity"
$ws.Range("B38").Value2 = 9
$ws.Range("D38").Value2 = "trauma - law enforcement"
$ws.Range("E38").Value2 = 989
$ws.Range("I38").Value2 = "bawtry - barberlife"
$ws.Range("J38").Value2 = 167
$ws.Range("L38").Value2 = "cumbria - fitness"
$ws.Range("M38").Value2 = 15339.1705661
$ws.Range("A39").Value2 = "mental health - art and (health)care"
$ws.Range("B39").Value2 = 8
$ws.Range("D39").Value2 = "psychology of medical technology - trauma"
$ws.Range("E39").Value2 = 989
$ws.Range("I39").Value2 = "barberlife - themanclub"
$ws.Range("J39").Value2 = 163
$ws.Range("L39").Value2 = "photooftheday - london"
$ws.Range("M39").Value2 = 15288.4480468
$ws.Range("A40").Value2 = "research question - story sharing"
$ws.Range("B40").Value2 = 8
$ws.Range("D40").Value2 = "cultural difference - map-making"
$ws.Range("E40").Value2 = 989
$ws.Range("I40").Value2 = "barberlife - barbershop"
$ws.Range("J40").Value2 = 163
$ws.Range("L40").Value2 = "london - yorkshire"
$ws.Range("M40").Value2 = 15167.2503872
$ws.Range("A41").Value2 = "methodology - community-based care"
$ws.Range("B41").Value2 = 8
$ws.Range("D41").Value2 = "homemade paint - hands-on/DIY"
$ws.Range("E41").Value2 = 989
$ws.Range("I41").Value2 = "barberlife - apothecary87"
$ws.Range("J41").Value2 = 154
$ws.Range("L41").Value2 = "london - homecinema"
$ws.Range("M41").Value2 = 14945.9185136
$ws.Range("A42").Value2 = "legality - safety"
$ws.Range("B42").Value2 = 8
$ws.Range("D42").Value2 = "moringa - blood regulatory function"
$ws.Range("E42").Value2 = 989
$ws.Range("I42").Value2 = "bawtry - tgr"
$ws.Range("J42").Value2 = 152
$ws.Range("L42").Value2 = "projector - london"
$ws.Range("M42").Value2 = 14945.9185136
$ws.Range("A43").Value2 = "crisis - resource strain"
$ws.Range("B43").Value2 = 8
$ws.Range("D43").Value2 = "tinkering is easier when the device is cheaper - hands-on/DIY"
$ws.Range("E43").Value2 = 989
$ws.Range("I43").Value2 = "thegentlemansretreat - tgr"
$ws.Range("J43").Value2 = 152
$ws.Range("L43").Value2 = "fitness - somerset"
$ws.Range("M43").Value2 = 14896.6659868
$ws.Range("A44").Value2 = "regulation - legality"
$ws.Range("B44").Value2 = 8
$ws.Range("D44").Value2 = "supporting not fixing - making rules for spaces"
$ws.Range("E44").Value2 = 989
$ws.Range("I44").Value2 = "themanclub - tgr"
$ws.Range("J44").Value2 = 151
$ws.Range("L44").Value2 = "beer - camping"
$ws.Range("M44").Value2 = 14511.9284747
$ws.Range("A45").Value2 = "mental health - suicide"
$ws.Range("B45").Value2 = 7
$ws.Range("D45").Value2 = "tolerance - social design"
$ws.Range("E45").Value2 = 989
$ws.Range("I45").Value2 = "tgr - barbershop"
$ws.Range("J45").Value2 = 151
$ws.Range("L45").Value2 = "humberston - cleethorpes"
$ws.Range("M45").Value2 = 14266

# --- Column widths to fit the long node/edge labels ---
$ws.Range("A:A").ColumnWidth = 35.753488372093
$ws.Range("B:C").ColumnWidth = 10.706976744186
$ws.Range("D:D").ColumnWidth = 51.4651162790698
$ws.Range("E:H").ColumnWidth = 10.706976744186
$ws.Range("I:I").ColumnWidth = 31.5767441860465
$ws.Range("J:K").ColumnWidth = 10.706976744186
$ws.Range("L:L").ColumnWidth = 26.893023255814

# --- Make "top20" the active sheet / selection, matching the authored workbook ---
$ws.Activate()
$ws.Range("A12").Select()

Write-Output "top20 sheet created with $($ws.UsedRange.Cells.Count) used cells"
